# Apply the data refresh described in the commit "Update countries & provincias Spain".
# The workbook is a single-sheet COVID-19 country dashboard ("Pais"). The update refreshes
# several countries' case counts; because the sheet is sorted by total cases, a handful of
# rows also swap which country they display.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Marzo de 2020 a las 07:25"

# --- Refresh per-country rows (columns: A=Pais, B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 142224
$ws.Cells.Item(4, 3).Value = 177
$ws.Cells.Item(4, 4).Value = 4559
$ws.Cells.Item(4, 5).Value = 135180
$ws.Cells.Item(4, 6).Value = 2970
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 2485

# Row 37: Tailandia
$ws.Cells.Item(37, 1).Value = "Tailandia"
$ws.Cells.Item(37, 2).Value = 1524
$ws.Cells.Item(37, 3).Value = 136
$ws.Cells.Item(37, 4).Value = 229
$ws.Cells.Item(37, 5).Value = 1288
$ws.Cells.Item(37, 6).Value = 11
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 7

# Row 38: Filipinas
$ws.Cells.Item(38, 1).Value = "Filipinas"
$ws.Cells.Item(38, 2).Value = 1418
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 42
$ws.Cells.Item(38, 5).Value = 1305
$ws.Cells.Item(38, 6).Value = 1
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 71

# Row 61: Nueva Zelanda
$ws.Cells.Item(61, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(61, 2).Value = 589
$ws.Cells.Item(61, 3).Value = 75
$ws.Cells.Item(61, 4).Value = 63
$ws.Cells.Item(61, 5).Value = 525
$ws.Cells.Item(61, 6).Value = 2
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 1

# Row 62: Emiratos Arabes Unidos
$ws.Cells.Item(62, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(62, 2).Value = 570
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 58
$ws.Cells.Item(62, 5).Value = 509
$ws.Cells.Item(62, 6).Value = 2
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 3

# Row 63: Irak
$ws.Cells.Item(63, 1).Value = "Irak"
$ws.Cells.Item(63, 2).Value = 547
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 143
$ws.Cells.Item(63, 5).Value = 362
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 42

# Row 69: Hungria
$ws.Cells.Item(69, 1).Value = "Hungria"
$ws.Cells.Item(69, 2).Value = 447
$ws.Cells.Item(69, 3).Value = 39
$ws.Cells.Item(69, 4).Value = 34
$ws.Cells.Item(69, 5).Value = 398
$ws.Cells.Item(69, 6).Value = 6
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 15

# Row 70: Libano
$ws.Cells.Item(70, 1).Value = "Libano"
$ws.Cells.Item(70, 2).Value = 438
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(70, 4).Value = 30
$ws.Cells.Item(70, 5).Value = 398
$ws.Cells.Item(70, 6).Value = 4
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 10

# Row 71: Armenia
$ws.Cells.Item(71, 1).Value = "Armenia"
$ws.Cells.Item(71, 2).Value = 424
$ws.Cells.Item(71, 3).Value = 0
$ws.Cells.Item(71, 4).Value = 30
$ws.Cells.Item(71, 5).Value = 391
$ws.Cells.Item(71, 6).Value = 6
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 3

# Row 91: Vietnam
$ws.Cells.Item(91, 1).Value = "Vietnam"
$ws.Cells.Item(91, 2).Value = 194
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 4).Value = 52
$ws.Cells.Item(91, 5).Value = 142
$ws.Cells.Item(91, 6).Value = 3
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 0

# Row 107: Mauricio
$ws.Cells.Item(107, 1).Value = "Mauricio"
$ws.Cells.Item(107, 2).Value = 110
$ws.Cells.Item(107, 3).Value = 3
$ws.Cells.Item(107, 4).Value = 0
$ws.Cells.Item(107, 5).Value = 107
$ws.Cells.Item(107, 6).Value = 1
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 3

# Row 108: Honduras
$ws.Cells.Item(108, 1).Value = "Honduras"
$ws.Cells.Item(108, 2).Value = 110
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 3
$ws.Cells.Item(108, 5).Value = 104
$ws.Cells.Item(108, 6).Value = 4
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 3

# Row 109: Estado de Palestina
$ws.Cells.Item(109, 1).Value = "Estado de Palestina"
$ws.Cells.Item(109, 2).Value = 109
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 18
$ws.Cells.Item(109, 5).Value = 90
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 1

# Row 112: Bolivia
$ws.Cells.Item(112, 1).Value = "Bolivia"
$ws.Cells.Item(112, 2).Value = 96
$ws.Cells.Item(112, 3).Value = 15
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(112, 5).Value = 95
$ws.Cells.Item(112, 6).Value = 3
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 1

# Row 113: Bielorrusia
$ws.Cells.Item(113, 1).Value = "Bielorrusia"
$ws.Cells.Item(113, 2).Value = 94
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 32
$ws.Cells.Item(113, 5).Value = 62
$ws.Cells.Item(113, 6).Value = 2
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 0

# Row 114: Martinica
$ws.Cells.Item(114, 1).Value = "Martinica"
$ws.Cells.Item(114, 2).Value = 93
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 92
$ws.Cells.Item(114, 6).Value = 12
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 1

# Row 115: Georgia
$ws.Cells.Item(115, 1).Value = "Georgia"
$ws.Cells.Item(115, 2).Value = 91
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 18
$ws.Cells.Item(115, 5).Value = 73
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 0

# Row 116: Montenegro
$ws.Cells.Item(116, 1).Value = "Montenegro"
$ws.Cells.Item(116, 2).Value = 85
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 84
$ws.Cells.Item(116, 6).Value = 1
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 1

# Row 117: Kirguistan
$ws.Cells.Item(117, 1).Value = "Kirguistan"
$ws.Cells.Item(117, 2).Value = 84
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 84
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 0

# Row 135: Barbados
$ws.Cells.Item(135, 1).Value = "Barbados"
$ws.Cells.Item(135, 2).Value = 33
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 0
$ws.Cells.Item(135, 5).Value = 33
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 0

# Row 136: Uganda
$ws.Cells.Item(136, 1).Value = "Uganda"
$ws.Cells.Item(136, 2).Value = 33
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 0
$ws.Cells.Item(136, 5).Value = 33
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 0
